$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 updates (add filter in simulation)
$ws.Range("B21").Value = "hispnditFilter"

$ws.Range("E21").Value = 0.75
$ws.Range("G21").Value = 1
$ws.Range("J21").Value = 0.5
$ws.Range("K21").Value = 1
$ws.Range("M21").Value = 0.3333333333333333
$ws.Range("O21").Value = 0.9795918367346941
$ws.Range("P21").Value = 0.9869281045751634
$ws.Range("Q21").Value = 0.7777777777777778
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 1
